$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rows 55 and 56 had their match data swapped (home/away teams, scores,
#    odds and match URL) -- dates/meta columns (A-E, I, K, M, O, Q, S, U)
#    stay in place.
# ---------------------------------------------------------------------------
$f55 = $ws.Range("F55").Value2
$g55 = $ws.Range("G55").Value2
$h55 = $ws.Range("H55").Value2
$j55 = $ws.Range("J55").Value2
$l55 = $ws.Range("L55").Value2
$n55 = $ws.Range("N55").Value2
$p55 = $ws.Range("P55").Value2
$r55 = $ws.Range("R55").Value2
$t55 = $ws.Range("T55").Value2
$v55 = $ws.Range("V55").Value2

$f56 = $ws.Range("F56").Value2
$g56 = $ws.Range("G56").Value2
$h56 = $ws.Range("H56").Value2
$j56 = $ws.Range("J56").Value2
$l56 = $ws.Range("L56").Value2
$n56 = $ws.Range("N56").Value2
$p56 = $ws.Range("P56").Value2
$r56 = $ws.Range("R56").Value2
$t56 = $ws.Range("T56").Value2
$v56 = $ws.Range("V56").Value2

$ws.Range("F55").Value = $f56
$ws.Range("G55").Value = $g56
$ws.Range("H55").Value = $h56
$ws.Range("J55").Value = $j56
$ws.Range("L55").Value = $l56
$ws.Range("N55").Value = $n56
$ws.Range("P55").Value = $p56
$ws.Range("R55").Value = $r56
$ws.Range("T55").Value = $t56
$ws.Range("V55").Value = $v56

$ws.Range("F56").Value = $f55
$ws.Range("G56").Value = $g55
$ws.Range("H56").Value = $h55
$ws.Range("J56").Value = $j55
$ws.Range("L56").Value = $l55
$ws.Range("N56").Value = $n55
$ws.Range("P56").Value = $p55
$ws.Range("R56").Value = $r55
$ws.Range("T56").Value = $t55
$ws.Range("V56").Value = $v55

# ---------------------------------------------------------------------------
# 2) A new match row (75) was appended at the bottom of the sheet. Copy the
#    formatting from the previous last row (74) so styles (bold/border on A,
#    date format on E) carry over, then overwrite with the new values.
# ---------------------------------------------------------------------------
$ws.Range("A74:V74").Copy($ws.Range("A75:V75"))

$ws.Range("A75").Value = 74
$ws.Range("B75").Value = "armenia"
$ws.Range("C75").Value = "premier-league"
$ws.Range("D75").Value = "2023-2024"
$ws.Range("E75").Value = 45238.54166666666
$ws.Range("F75").Value = "Shirak Gyumri"
$ws.Range("G75").Value = 1
$ws.Range("H75").Value = "Ararat Yerevan"
$ws.Range("I75").Value = 2
$ws.Range("J75").Value = 2.7
$ws.Range("K75").Value = "07/11/2023 01:12"
$ws.Range("L75").Value = 2.6
$ws.Range("M75").Value = "08/11/2023 12:55"
$ws.Range("N75").Value = 3
$ws.Range("O75").Value = "07/11/2023 01:12"
$ws.Range("P75").Value = 3.17
$ws.Range("Q75").Value = "08/11/2023 12:59"
$ws.Range("R75").Value = 2.51
$ws.Range("S75").Value = "07/11/2023 01:12"
$ws.Range("T75").Value = 2.8
$ws.Range("U75").Value = "08/11/2023 12:54"
$ws.Range("V75").Value = "https://www.betexplorer.com/football/armenia/premier-league/shirak-gyumri-ararat-yerevan/WMclsZ4n/"
